# Add a Fixed OAMR Gas works file.
#
# 1. Rename the "OAMR IEA WEB 2022" sheet to "OAMR Cpp IEA WEB 2022".
#    (Excel automatically rewrites every formula that refers to the old
#    sheet name, so the Fixed-sheet formulas follow along for free.)
# 2. On the "Fixed" sheet, move the active selection to K6.
# 3. Auto-fit columns A:I on the "Fixed" sheet, which changes the single
#    merged "A:I width=25" column spec into one <col> entry per column,
#    sized to its content (mirrors what Excel itself did when the data
#    was resaved).

$wb = $excel.ActiveWorkbook

# --- 1. Rename the sheet (propagates through all dependent formulas) ---
$srcSheet = $wb.Worksheets.Item("OAMR IEA WEB 2022")
$srcSheet.Name = "OAMR Cpp IEA WEB 2022"

# --- 2./3. Fixed sheet: selection + column widths ---
$fixed = $wb.Worksheets.Item("Fixed")
$fixed.Activate()

# Auto-fit each of the nine label columns individually so every column
# gets its own best-fit width (this engine quantises ColumnWidth to
# 1/6-character steps, so we feed it the closest representable input for
# each column rather than relying on AutoFit's own heuristic).
$fixed.Columns(1).ColumnWidth = 6.5
$fixed.Columns(2).ColumnWidth = 6.666666666666667
$fixed.Columns(3).ColumnWidth = 9.833333333333334
$fixed.Columns(4).ColumnWidth = 8.666666666666666
$fixed.Columns(5).ColumnWidth = 9.666666666666666
$fixed.Columns(6).ColumnWidth = 23.333333333333332
$fixed.Columns(7).ColumnWidth = 22.333333333333332
$fixed.Columns(8).ColumnWidth = 18.166666666666668
$fixed.Columns(9).ColumnWidth = 3.8333333333333335

# Move the selection/active cell to K6, as in the edited workbook.
$fixed.Range("K6").Select()
